$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "18Ballard" data currently in row 3 (row 4 "20Ballard" stays put,
# leaving a gap at row 3) and append a fresh "18Ballard" row at row 5.
$ws.Range("A3:G3").ClearContents()

$ws.Range("A5").Value = "18Ballard"
$ws.Range("B5").Value = 18
$ws.Range("C5").Value = 145
$ws.Range("D5").Value = 0.36
$ws.Range("E5").Value = 0.45
$ws.Range("F5").Value = 0.04
$ws.Range("G5").Value = 1.02

# Match the workbook's final selection state (whole row 5 selected).
$ws.Rows.Item(5).Select()
